$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(4544,0,5.88,13104,35858.28,8.15,53520.310000000005)
    3 = @(4736,39.36,43.65,10224,35468.67,0,50511.67999999999)
    4 = @(4862,871.79,919.58,10800,35693.31,81.67,53228.34999999999)
    5 = @(7236,10591.85,12151.5,14256,42595.48,2431.5,89262.33)
    6 = @(18398,52350.26,62941.53,51984,101223.82,24516.14,311413.75)
    7 = @(37016,114924.91,131737.7,149472,296065.26,119522.57,848738.44)
    8 = @(64970,254796.6,246824.99,310896,577702.4,287319.94,1742509.93)
    9 = @(71652,366853.63,317068.5,232704,430127.94,175900.63,1594306.7000000002)
    10 = @(18942,83955.49,80286.24,67680,108885.49,41099.06,400848.28)
    11 = @(9676,34554.29,34860.74,42048,63945.23,25182.49,210266.75)
    12 = @(7524,24203.62,26909.73,26928,45729.61,15694.42,146989.38000000003)
    13 = @(11090,40971.58,50641.16,23760,41391.22,8783.19,176637.15000000002)
    14 = @(20524,71516.91,113937.47,31824,54037.38,4531.6,296371.36)
    15 = @(28966,59819.23,155077.26,40896,96901.17,2208.43,383868.08999999997)
    16 = @(35200,66494.81,154748.2,75456,183380.68,1340.87,516620.56)
    17 = @(46912,90214.77,172740.33,129168,309365.09,532.26,748932.45)
    18 = @(57282,118517.95,153137.12,175104,452247.29,105.84,956394.2)
    19 = @(64190,140877.81,120788.95,232992,589716.82,9.65,1148575.23)
    20 = @(74316,234079.1,142299.29,259920,663277.97,27.62,1373919.98)
    21 = @(70850,191141.86,123446.99,251856,655664.5,0,1292959.35)
    22 = @(61274,130059.08,116683.37,212544,573511.06,0,1094071.51)
    23 = @(51730,103329.06,152569.41,174528,433865.24,0,916021.71)
    24 = @(37356,71781.4,131021.31,105120,264276.29,0,609555)
    25 = @(23744,43613.81,98595.26,59328,136454.52,0,361735.58999999997)
    26 = @(832990,2305559.1700000004,2599436.16,2702592,6227384.72,709296.0299999999,15377258.08)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $col = $i + 2
        $ws.Cells.Item($r, $col).Value = $rowVals[$i]
    }
}

Write-Host "Updated rows 2-26, columns B-H"